# DaySale report update: replace "ABIMOL 300MG 5 RECTAL SUPP." row with
# "ALBENDAZOLE 400MG 6 TABS" (shifting the following rows up) and insert a
# new "DIAMICRON 60MG M.R. 30 SCORED TAB" row after "DEPAKINE CHRONO..."
# (before "EPICOPRED..."). Only rows 7-13 actually change text; everything
# from row 14 downward is untouched. Update the running total (P33) and the
# generated-at timestamp (A34) to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7: was ABIMOL -> now ALBENDAZOLE (previously row 8's data)
$ws.Range("C7").Value = "ALBENDAZOLE 400MG 6 TABS"
$ws.Range("H7").Value = "0:0"
$ws.Range("N7").Value = "46.00"
$ws.Range("P7").Value = "92.0000"
$ws.Range("Q7").Value = "2:0"

# Row 8: was ALBENDAZOLE -> now ANTINAL (previously row 9's data)
$ws.Range("C8").Value = "ANTINAL 220MG/5ML 60ML SUSP."
$ws.Range("H8").Value = "3:0"
$ws.Range("N8").Value = "24.00"
$ws.Range("P8").Value = "24.0000"
$ws.Range("Q8").Value = "1:0"

# Row 9: was ANTINAL -> now COGINTOL (previously row 10's data)
$ws.Range("C9").Value = "COGINTOL 20 TAB."
$ws.Range("H9").Value = "0:0"
$ws.Range("N9").Value = "40.00"
$ws.Range("P9").Value = "40.0000"

# Row 10: was COGINTOL -> now COLONA (previously row 11's data)
$ws.Range("C10").Value = "COLONA 30 F.C.TAB"
$ws.Range("H10").Value = "0:1"
$ws.Range("N10").Value = "69.00"
$ws.Range("P10").Value = "69.0000"

# Row 11: was COLONA -> now DENSITIN (previously row 12's data)
$ws.Range("C11").Value = "DENSITIN 30 CAPS"
$ws.Range("H11").Value = "0:0"
$ws.Range("N11").Value = "96.00"
$ws.Range("P11").Value = "63.3600"
$ws.Range("Q11").Value = "0:2"

# Row 12: was DENSITIN -> now DEPAKINE CHRONO (previously row 13's data)
$ws.Range("C12").Value = "DEPAKINE CHRONO 500MG 30 SCORED PROLONGED REL. F.C. TAB."
$ws.Range("N12").Value = "144.00"
$ws.Range("P12").Value = "144.0000"
$ws.Range("Q12").Value = "1:0"

# Row 13: was DEPAKINE CHRONO -> now the new DIAMICRON row
$ws.Range("C13").Value = "DIAMICRON 60MG M.R. 30 SCORED TAB"
$ws.Range("H13").Value = "1:2"
$ws.Range("N13").Value = "156.00"
$ws.Range("P13").Value = "156.0000"

# Rows 14-32 (EPICOPRED ... محلول ملح) are unchanged.

# Running total of column P (price) shifts from 1596.73 to 1737.73
# (-15.00 for the removed ABIMOL row, +156.00 for the new DIAMICRON row).
$ws.Range("P33").Value = 1737.73

# Footer generation timestamp.
$ws.Range("A34").Value = "Monday, 15 September, 2025 10:34 AM"
